$d = $word.ActiveDocument

# --- Step 1: Insert the new "Act 2" section paragraphs ------------------
# Locate the empty paragraph that sits right after the big intro story
# paragraph ("...Maybe the others awoke before you...") and right before
# "Freight Bay door:". We insert six new paragraphs after it:
#   '', 'Act 2', '', 'Environmental Controls keypad:', '', ''
# each inheriting the en-FI / left-aligned formatting already used there.

$introPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Maybe the others awoke before you*") {
        $introPara = $i
    }
}

if ($introPara -eq $null) {
    throw "Could not locate intro story paragraph"
}

$idx = $introPara + 1

$newParas = @("", "Act 2", "", "Environmental Controls keypad:", "", "")
foreach ($t in $newParas) {
    $r = $d.Paragraphs($idx).Range
    $r.InsertParagraphAfter()
    $idx = $idx + 1
    if ($t.Length -gt 0) {
        $d.Paragraphs($idx).Range.InsertAfter($t)
    }
}

# --- Step 2: Remove the trailing empty paragraphs after the Freight Bay
# door description paragraph -------------------------------------------

$doorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*without a careful approach*") {
        $doorPara = $i
    }
}

if ($doorPara -eq $null) {
    throw "Could not locate freight bay door paragraph"
}

if ($doorPara -lt $d.Paragraphs.Count) {
    $startPara = $d.Paragraphs($doorPara + 1)
    $endPara = $d.Paragraphs($d.Paragraphs.Count)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
